$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("E2").Value = 1156.8090495
$wsSchedule.Range("F2").Value = 25.50284500661376
$wsSchedule.Range("E3").Value = 309.4335659999999
$wsSchedule.Range("F3").Value = 20.46518293650793
$wsSchedule.Range("A4").Value = 46039.3125
$wsSchedule.Range("B4").Value = 46039.8125
$wsSchedule.Range("E4").Value = -117.1090635
$wsSchedule.Range("F4").Value = -2.581769477513227

# --- Detailed sheet updates ---
$wsDetailed.Range("B33").Value = 25.59822
$wsDetailed.Range("B34").Value = 18.17021
$wsDetailed.Range("B35").Value = -34.28578
$wsDetailed.Range("C35").Value = "historical"
$wsDetailed.Range("B36").Value = -9.77904
$wsDetailed.Range("C36").Value = "historical"
$wsDetailed.Range("B37").Value = -6.84011
$wsDetailed.Range("B38").Value = -3.13143
$wsDetailed.Range("B39").Value = -3.05909
$wsDetailed.Range("B40").Value = -3.12465
$wsDetailed.Range("B41").Value = 5.07885
$wsDetailed.Range("B42").Value = 10.55942
$wsDetailed.Range("B43").Value = 21.58386
$wsDetailed.Range("B44").Value = 0.19633
$wsDetailed.Range("B45").Value = 57.04922
$wsDetailed.Range("B47").Value = 57.03893
$wsDetailed.Range("B48").Value = 56.98
$wsDetailed.Range("B59").Value = 56.98
$wsDetailed.Range("B62").Value = 57.06018
$wsDetailed.Range("B65").Value = 0.51
$wsDetailed.Range("E65").Value = "ON"
$wsDetailed.Range("B66").Value = -5.51
$wsDetailed.Range("B67").Value = -4.81333
$wsDetailed.Range("B68").Value = 0.51003
$wsDetailed.Range("B69").Value = 0
$wsDetailed.Range("B70").Value = -0.87926
$wsDetailed.Range("B71").Value = -0.87869
$wsDetailed.Range("B72").Value = -0.85809
$wsDetailed.Range("B73").Value = -0.84198
$wsDetailed.Range("B74").Value = -3.6481
$wsDetailed.Range("B75").Value = -4.81333
$wsDetailed.Range("B76").Value = -8.5
$wsDetailed.Range("B77").Value = -10
$wsDetailed.Range("B78").Value = -14
$wsDetailed.Range("B79").Value = -14
$wsDetailed.Range("B80").Value = -12.11173
$wsDetailed.Range("B81").Value = -7.92844
$wsDetailed.Range("B82").Value = -2.88008
$wsDetailed.Range("B83").Value = -6.76678
$wsDetailed.Range("B84").Value = -7.78482
$wsDetailed.Range("B85").Value = -4.96609
$wsDetailed.Range("B86").Value = -6.19024
$wsDetailed.Range("B87").Value = -6
$wsDetailed.Range("B88").Value = 2.23907
$wsDetailed.Range("B89").Value = 17.98199
$wsDetailed.Range("E89").Value = "OFF"
$wsDetailed.Range("B90").Value = 32.58868
$wsDetailed.Range("B91").Value = 29.51425
$wsDetailed.Range("B92").Value = 29.67769
$wsDetailed.Range("B95").Value = 64.8901
$wsDetailed.Range("B96").Value = 64.8901
$wsDetailed.Range("B97").Value = 64.8901
